$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency prices and 1h volume percentages (scraped refresh)
# Values that parse as plain numbers get a leading apostrophe so Excel
# keeps them as text (matching the source data's text-formatted cells).
$ws.Range("D2").Value = "26.900.55"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.844.61"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'309.63"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "'0.4694"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").Value = "'0.07153"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'0.9263"
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D11").Value = "'19.57"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "'0.07692"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "1.914.00"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("D14").Value = "'5.289"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'6.397"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "'88.24"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'0.000008633"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "26.930.95"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "'5.015"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D25").Value = "'151.92"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'18.26"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").Value = "'2.007"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "'114.19"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'4.880"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").Value = "'3.219"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").Value = "'1.180"
$ws.Range("E32").Value = "  +6.30%  "
$ws.Range("D33").Value = "'0.7462"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'2.786"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").Value = "'4.471"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").Value = "'1.085"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").Value = "'0.01936"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'0.05207"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").Value = "'2.957"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").Value = "'0.5204"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("D41").Value = "'6.975"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("D42").Value = "'0.1513"
$ws.Range("D43").Value = "'8.158"
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("D44").Value = "'10.49"
$ws.Range("E44").Value = "  +5.63%  "
$ws.Range("D45").Value = "'0.4695"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'1.007"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "'101.22"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Value = "'1.596"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "'65.40"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "'0.06036"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "'0.8948"
$ws.Range("E51").Value = "  +5.69%  "
